$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P holds 2022 data, one row per metric (mirrors the existing
# 2010-2021 columns D:O). Copy column O's formatting into P first so each
# new cell inherits the same number format / font / fill as its row, then
# overwrite the value with the 2022 figure.
$rows = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17)
foreach ($r in $rows) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 11.4
$ws.Range("P6").Value = 12.6
$ws.Range("P7").Value = 9.8
$ws.Range("P8").Value = 11.4
$ws.Range("P9").Value = 5.4
$ws.Range("P10").Value = 4.7
$ws.Range("P11").Value = 3.4
$ws.Range("P12").Value = 17.7
$ws.Range("P13").Value = 20.5
$ws.Range("P14").Value = 8.4
$ws.Range("P16").Value = 12.9
$ws.Range("P17").Value = 10.5

# Restore the selection the author left behind.
$ws.Range("Q4").Select()
